# daily auto push: 2026-01-25 18:44 UTC
# Two new 3-hourly readings were recorded for 2026/01/25 (Sun) and
# 2026/01/26 (Mon) just before the 2026/12/29 block; insert them at the
# top of the existing "tail" data block (row 713) which pushes the rest
# of the table down by two rows (A1:D754 -> A1:D756).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 713 (each Insert() shifts row 713-on down
# by one, so calling it twice makes room for both new rows).
$ws.Rows.Item(713).Insert()
$ws.Rows.Item(713).Insert()

function Set-TextCell($row, $col, $text) {
    # Column A holds yyyy/mm/dd-looking strings that must stay literal
    # text (matching the rest of the sheet) instead of being silently
    # auto-converted into a date serial by the normal .Value setter.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# New row 713: 2026/01/25 (Sun), hour 23
Set-TextCell 713 1 "2026/01/25"
$ws.Cells.Item(713, 2).Value = "日"
$ws.Cells.Item(713, 3).Value = 23
$ws.Cells.Item(713, 4).Value = 201

# New row 714: 2026/01/26 (Mon), hour 2
Set-TextCell 714 1 "2026/01/26"
$ws.Cells.Item(714, 2).Value = "月"
$ws.Cells.Item(714, 3).Value = 2
$ws.Cells.Item(714, 4).Value = 201
